# "working on 2001 - 2300"
# Fill in the missing "paid" amount for the 3201-3500 pay-line row (row 15),
# mark it + the next row (3501-3800) as moved-to-the-"T" file ("yes" in L),
# and give row 16 its own (non-shared) formulas since the 300 -> 390
# increment no longer matches the shared H9:H16 pattern. Then extend the
# running-total SUM down into row 17 and move the selection/view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15: the 280 payment that was missing.
$ws.Range("I15").Value = 280

# Row 15 now belongs to the "moved to file" ("yes") bucket.
$ws.Range("L15").Value = "yes"

# Row 16: the increment changed from +300 to +390, so it can no longer
# share row 9's formula - give it its own literal formula (H9:H15 keep
# their original G{row}+300-1 formula/value, just no longer grouped
# with H16 as one shared formula).
$ws.Range("H16").Formula = "=G16+390-1"

# Row 16's paid amount, pro-rated from the 300 -> 390 change.
$ws.Range("I16").Formula = "=280*390/300"

# Row 16 is also moved to the "T" file.
$ws.Range("L16").Value = "yes"

# Carry the running total (SUM(I7:I16)) down onto row 17.
$ws.Range("I17").Formula = "=SUM(I7:I16)"

# Move the view/selection to reflect where work continued.
$ws.Range("O20").Select()
